$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 11 values (recomputed baseline numbers)
$ws.Range("C11").Value = 14890.763749387301
$ws.Range("D11").Value = 22335.029821411099
$ws.Range("E11").Value = 1279.6239764917827
$ws.Range("F11").Value = 17112.929523077662

# Fill in row 12 with a new PPO2 / PRIVE+VOL entry
$ws.Range("A12").Value = "PPO2"
$ws.Range("B12").Value = "PRIVE+VOL"
$ws.Range("C12").Value = 13963.501355182199
$ws.Range("D12").Value = 21294.655943123598
$ws.Range("E12").Value = 1295.9528501222114
$ws.Range("F12").Value = 16826.283226245319

# Update the sheet selection to match the author's final state
$ws.Range("C11:F11").Select()
